$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040809482733422
$ws.Range("D2").Value = 1.04365670735475
$ws.Range("E2").Value = 1.044466440537053
$ws.Range("F2").Value = 1.049227967760451
$ws.Range("I2").Value = 1.040980793633511
$ws.Range("J2").Value = 1.045894083157781
$ws.Range("K2").Value = 1.046429957196339
$ws.Range("L2").Value = 1.047237410918254
$ws.Range("M2").Value = 1.05198561125979
$ws.Range("N2").Value = 1.019101556494424

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041960692671536
$ws.Range("D3").Value = 1.044552503137746
$ws.Range("E3").Value = 1.045566763666612
$ws.Range("F3").Value = 1.050587390237693
$ws.Range("I3").Value = 1.041334458434953
$ws.Range("J3").Value = 1.046690057353124
$ws.Range("K3").Value = 1.047136610946707
$ws.Range("L3").Value = 1.048148223623396
$ws.Range("M3").Value = 1.053155824479918
$ws.Range("N3").Value = 1.019374910352988

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042705259563075
$ws.Range("D4").Value = 1.045131795468497
$ws.Range("E4").Value = 1.046278753314007
$ws.Range("F4").Value = 1.051467046761273
$ws.Range("I4").Value = 1.041561935351102
$ws.Range("J4").Value = 1.047204229612658
$ws.Range("K4").Value = 1.04759288404518
$ws.Range("L4").Value = 1.048736999023224
$ws.Range("M4").Value = 1.053912516035989
$ws.Range("N4").Value = 1.019551279211235

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043018195191657
$ws.Range("D5").Value = 1.045375247505963
$ws.Range("E5").Value = 1.04657807665914
$ws.Range("F5").Value = 1.051836861748585
$ws.Range("I5").Value = 1.041657239700608
$ws.Range("J5").Value = 1.047420178904262
$ws.Range("K5").Value = 1.047784467679103
$ws.Range("L5").Value = 1.04898438215293
$ws.Range("M5").Value = 1.054230508322449
$ws.Range("N5").Value = 1.01962530295428

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043070733804669
$ws.Range("D6").Value = 1.045416119359757
$ws.Range("E6").Value = 1.046628334571213
$ws.Range("F6").Value = 1.051898955825761
$ws.Range("I6").Value = 1.04167322256315
$ws.Range("J6").Value = 1.047456425531355
$ws.Range("K6").Value = 1.047816621741263
$ws.Range("L6").Value = 1.049025910806302
$ws.Range("M6").Value = 1.054283893608947
$ws.Range("N6").Value = 1.019637724734773

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042709441337496
$ws.Range("D7").Value = 1.045135048809645
$ws.Range("E7").Value = 1.046282752874282
$ws.Range("F7").Value = 1.051471988214749
$ws.Range("I7").Value = 1.041563210096506
$ws.Range("J7").Value = 1.047207115956497
$ws.Range("K7").Value = 1.047595444911837
$ws.Range("L7").Value = 1.048740305109745
$ws.Range("M7").Value = 1.053916765536362
$ws.Range("N7").Value = 1.019552268798136

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041198611610893
$ws.Range("D8").Value = 1.043959518108343
$ws.Range("E8").Value = 1.044838299175187
$ws.Range("F8").Value = 1.049687388626935
$ws.Range("I8").Value = 1.041100599848025
$ws.Range("J8").Value = 1.046163268437026
$ws.Range("K8").Value = 1.046668977375381
$ws.Range("L8").Value = 1.047545345379783
$ws.Range("M8").Value = 1.052381198280875
$ws.Range("N8").Value = 1.019194043310705

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03853362923548
$ws.Range("D9").Value = 1.041885389020909
$ws.Range("E9").Value = 1.042292985751576
$ws.Range("F9").Value = 1.046542719660442
$ws.Range("I9").Value = 1.040274921699988
$ws.Range("J9").Value = 1.044317116937615
$ws.Range("K9").Value = 1.045028886677315
$ws.Range("L9").Value = 1.045435170737105
$ws.Range("M9").Value = 1.049671284974589
$ws.Range("N9").Value = 1.018558893062901

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036755038419292
$ws.Range("D10").Value = 1.040500776561428
$ws.Range("E10").Value = 1.040596029447592
$ws.Range("F10").Value = 1.044446113965003
$ws.Range("I10").Value = 1.039717375823165
$ws.Range("J10").Value = 1.043081736290333
$ws.Range("K10").Value = 1.043930371090301
$ws.Range("L10").Value = 1.044025286921705
$ws.Range("M10").Value = 1.047861804525329
$ws.Range("N10").Value = 1.018132814860289

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035984403326496
$ws.Range("D11").Value = 1.039900771323511
$ws.Range("E11").Value = 1.039861190927196
$ws.Range("F11").Value = 1.043538182988948
$ws.Range("I11").Value = 1.039474261695813
$ws.Range("J11").Value = 1.042545692711184
$ws.Range("K11").Value = 1.043453474053263
$ws.Range("L11").Value = 1.04341403917949
$ws.Range("M11").Value = 1.047077564372267
$ws.Range("N11").Value = 1.017947687329048

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035698078615486
$ws.Range("D12").Value = 1.039677832613339
$ws.Range("E12").Value = 1.039588230685654
$ws.Range("F12").Value = 1.043200920441951
$ws.Range("I12").Value = 1.039383703158199
$ws.Range("J12").Value = 1.042346413182723
$ws.Range("K12").Value = 1.043276146888323
$ws.Range("L12").Value = 1.04318687899601
$ws.Range("M12").Value = 1.046786150989669
$ws.Range("N12").Value = 1.017878827282352

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035759499739112
$ws.Range("D13").Value = 1.03972565691115
$ws.Range("E13").Value = 1.039646781995889
$ws.Range("F13").Value = 1.043273265208754
$ws.Range("I13").Value = 1.039403139829134
$ws.Range("J13").Value = 1.042389166992726
$ws.Range("K13").Value = 1.043314192620128
$ws.Range("L13").Value = 1.043235610880807
$ws.Range("M13").Value = 1.046848665251439
$ws.Range("N13").Value = 1.017893602320721

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035960737211423
$ws.Range("D14").Value = 1.039882344575583
$ws.Range("E14").Value = 1.039838628131655
$ws.Range("F14").Value = 1.043510305119661
$ws.Range("I14").Value = 1.039466781307737
$ws.Range("J14").Value = 1.042529223660804
$ws.Range("K14").Value = 1.043438819943353
$ws.Range("L14").Value = 1.043395264414234
$ws.Range("M14").Value = 1.047053478346724
$ws.Range("N14").Value = 1.017941997285116

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036084716048266
$ws.Range("D15").Value = 1.039978875701593
$ws.Range("E15").Value = 1.039956829690829
$ws.Range("F15").Value = 1.043656350933142
$ws.Range("I15").Value = 1.039505959093712
$ws.Range("J15").Value = 1.042615494748254
$ws.Range("K15").Value = 1.043515582218498
$ws.Range("L15").Value = 1.043493616866233
$ws.Range("M15").Value = 1.047179655566336
$ws.Range("N15").Value = 1.017971802362564

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036806172331435
$ws.Range("D16").Value = 1.040540587170291
$ws.Range("E16").Value = 1.040644797161968
$ws.Range("F16").Value = 1.044506368224369
$ws.Range("I16").Value = 1.039733474771023
$ws.Range("J16").Value = 1.04311728810749
$ws.Range("K16").Value = 1.043961995105838
$ws.Range("L16").Value = 1.044065837302058
$ws.Range("M16").Value = 1.047913836467141
$ws.Range("N16").Value = 1.018145087800661

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037258588988956
$ws.Range("D17").Value = 1.04089281029537
$ws.Range("E17").Value = 1.041076327671172
$ws.Range("F17").Value = 1.045039535884545
$ws.Range("I17").Value = 1.039875735522795
$ws.Range("J17").Value = 1.043431750110434
$ws.Range("K17").Value = 1.04424168728772
$ws.Range("L17").Value = 1.044424571824523
$ws.Range("M17").Value = 1.048374173003835
$ws.Range("N17").Value = 1.018253615491346

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037522428512127
$ws.Range("D18").Value = 1.041098211881798
$ws.Range("E18").Value = 1.041328027865413
$ws.Range("F18").Value = 1.04535051541372
$ws.Range("I18").Value = 1.0399585505055
$ws.Range("J18").Value = 1.043615062935441
$ws.Range("K18").Value = 1.04440470816384
$ws.Range("L18").Value = 1.044633742499358
$ws.Range("M18").Value = 1.04864260994438
$ws.Range("N18").Value = 1.018316856819821

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037612382993663
$ws.Range("D19").Value = 1.041168241050545
$ws.Range("E19").Value = 1.041413850457686
$ws.Range("F19").Value = 1.045456550121423
$ws.Range("I19").Value = 1.039986760597848
$ws.Range("J19").Value = 1.043677549675298
$ws.Range("K19").Value = 1.044460273948484
$ws.Range("L19").Value = 1.044705051941968
$ws.Range("M19").Value = 1.048734128295094
$ws.Range("N19").Value = 1.018338410136411

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037210053888993
$ws.Range("D20").Value = 1.040855024667378
$ws.Range("E20").Value = 1.041030028983201
$ws.Range("F20").Value = 1.044982332908004
$ws.Range("I20").Value = 1.039860489187442
$ws.Range("J20").Value = 1.043398022447418
$ws.Range("K20").Value = 1.044211691264373
$ws.Range("L20").Value = 1.044386090577714
$ws.Range("M20").Value = 1.048324790467764
$ws.Range("N20").Value = 1.018241977813132

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03590147992313
$ws.Range("D21").Value = 1.039836205906358
$ws.Range("E21").Value = 1.039782134474534
$ws.Range("F21").Value = 1.043440503259207
$ws.Range("I21").Value = 1.039448047527348
$ws.Range("J21").Value = 1.042487985142018
$ws.Range("K21").Value = 1.043402125457321
$ws.Range("L21").Value = 1.043348253630401
$ws.Range("M21").Value = 1.046993169101983
$ws.Range("N21").Value = 1.01792774881179

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035078283601921
$ws.Range("D22").Value = 1.039195228277595
$ws.Range("E22").Value = 1.038997483664521
$ws.Range("F22").Value = 1.042470996229415
$ws.Range("I22").Value = 1.039187252627743
$ws.Range("J22").Value = 1.04191482961727
$ws.Range("K22").Value = 1.042892040418448
$ws.Range("L22").Value = 1.042695054999145
$ws.Range("M22").Value = 1.04615527829909
$ws.Range("N22").Value = 1.017729628207828

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035514718422516
$ws.Range("D23").Value = 1.039535061485611
$ws.Range("E23").Value = 1.03941344721413
$ws.Range("F23").Value = 1.042984960667746
$ws.Range("I23").Value = 1.039325645143807
$ws.Range("J23").Value = 1.042218763490667
$ws.Range("K23").Value = 1.043162548815505
$ws.Range("L23").Value = 1.043041391922942
$ws.Range("M23").Value = 1.046599522463397
$ws.Range("N23").Value = 1.017834708166463

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037231984957527
$ws.Range("D24").Value = 1.040872098502438
$ws.Range("E24").Value = 1.041050949380053
$ws.Range("F24").Value = 1.04500818049515
$ws.Range("I24").Value = 1.039867378854484
$ws.Range("J24").Value = 1.043413262858417
$ws.Range("K24").Value = 1.044225245543008
$ws.Range("L24").Value = 1.044403478822179
$ws.Range("M24").Value = 1.048347104524386
$ws.Range("N24").Value = 1.018247236567455

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039222923782161
$ws.Range("D25").Value = 1.042421925313266
$ws.Range("E25").Value = 1.042951018640506
$ws.Range("F25").Value = 1.04735570723026
$ws.Range("I25").Value = 1.040489627097521
$ws.Range("J25").Value = 1.044795199365464
$ws.Range("K25").Value = 1.045453788162543
$ws.Range("L25").Value = 1.045981242601018
$ws.Range("M25").Value = 1.050372357628194
$ws.Range("N25").Value = 1.018723559660655
